# Actualización automática 2025-06-25 16:00:10
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M19").Value = 45.36
$wsVentasGrupo.Range("M29").Value = "2 de 27"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F19").Value = 411.7
$wsVentaMensual.Range("F29").Value = 13931.95

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 3372.76
$wsCumplimiento.Range("E16").Value = 9688.82
$wsCumplimiento.Range("F16").Value = 0.2582199090768498

$wsCumplimiento.Range("D19").Value = 20878.78
$wsCumplimiento.Range("E19").Value = 2621.220930050386
$wsCumplimiento.Range("F19").Value = 0.8884586882420704
